$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '26.407.73'
$ws.Cells.Item(2, 5).Value2 = '  -1.89%  '
$ws.Cells.Item(3, 4).Value2 = '1.792.59'
$ws.Cells.Item(3, 5).Value2 = '  -2.31%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value2 = '1.008'
$ws.Cells.Item(4, 5).Value2 = '  +0.32%  '
$ws.Cells.Item(5, 5).Value2 = '  +0.35%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = '307.01'
$ws.Cells.Item(6, 5).Value2 = '  -1.15%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = '0.4566'
$ws.Cells.Item(7, 5).Value2 = '  -1.01%  '
$ws.Cells.Item(8, 5).Value2 = '  -1.00%  '
$ws.Cells.Item(9, 2).Value2 = 'OKB'
$ws.Cells.Item(9, 3).Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = '46.82'
$ws.Cells.Item(9, 5).Value2 = '  +1.95%  '
$ws.Cells.Item(10, 2).Value2 = 'Dogecoin'
$ws.Cells.Item(10, 3).Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = '0.07090'
$ws.Cells.Item(10, 5).Value2 = '  -1.15%  '
$ws.Cells.Item(11, 2).Value2 = 'Polygon'
$ws.Cells.Item(11, 3).Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '0.8772'
$ws.Cells.Item(11, 5).Value2 = '  -0.64%  '
$ws.Cells.Item(12, 2).Value2 = 'TRON'
$ws.Cells.Item(12, 3).Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = '0.07884'
$ws.Cells.Item(12, 5).Value2 = '  +0.63%  '
$ws.Cells.Item(13, 2).Value2 = 'Solana'
$ws.Cells.Item(13, 3).Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = '19.53'
$ws.Cells.Item(13, 5).Value2 = '  -0.70%  '
$ws.Cells.Item(14, 2).Value2 = 'WrappedEther'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value2 = '1.801.51'
$ws.Cells.Item(14, 5).Value2 = '  -1.60%  '
$ws.Cells.Item(15, 2).Value2 = 'Polkadot'
$ws.Cells.Item(15, 3).Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = '5.275'
$ws.Cells.Item(15, 5).Value2 = '  -1.11%  '
$ws.Cells.Item(16, 2).Value2 = 'Chainlink'
$ws.Cells.Item(16, 3).Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '6.326'
$ws.Cells.Item(16, 5).Value2 = '  -1.02%  '
$ws.Cells.Item(17, 2).Value2 = 'Litecoin'
$ws.Cells.Item(17, 3).Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = '84.97'
$ws.Cells.Item(17, 5).Value2 = '  -4.10%  '
$ws.Cells.Item(18, 2).Value2 = 'BinanceUSD'
$ws.Cells.Item(18, 3).Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = '1.009'
$ws.Cells.Item(18, 5).Value2 = '  +0.40%  '
$ws.Cells.Item(19, 2).Value2 = 'ShibaInu'
$ws.Cells.Item(19, 3).Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = '0.000008522'
$ws.Cells.Item(19, 5).Value2 = '  -2.80%  '
$ws.Cells.Item(20, 2).Value2 = 'Dai'
$ws.Cells.Item(20, 3).Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = '1.008'
$ws.Cells.Item(20, 5).Value2 = '  +0.34%  '
$ws.Cells.Item(21, 2).Value2 = 'WrappedBTC'
$ws.Cells.Item(21, 3).Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(21, 4).Value2 = '26.434.37'
$ws.Cells.Item(21, 5).Value2 = '  -1.85%  '
$ws.Cells.Item(22, 2).Value2 = 'Avalanche'
$ws.Cells.Item(22, 3).Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '14.25'
$ws.Cells.Item(22, 5).Value2 = '  -1.78%  '
$ws.Cells.Item(23, 2).Value2 = 'Uniswap'
$ws.Cells.Item(23, 3).Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = '4.980'
$ws.Cells.Item(23, 5).Value2 = '  -0.77%  '
$ws.Cells.Item(24, 2).Value2 = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(24, 3).Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(24, 4).Value2 = '2.014.36'
$ws.Cells.Item(24, 5).Value2 = '  -2.03%  '
$ws.Cells.Item(25, 2).Value2 = 'Cosmos'
$ws.Cells.Item(25, 3).Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '10.49'
$ws.Cells.Item(25, 5).Value2 = '  +0.41%  '
$ws.Cells.Item(26, 2).Value2 = 'Toncoin'
$ws.Cells.Item(26, 3).Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = '1.985'
$ws.Cells.Item(26, 5).Value2 = '  +0.59%  '
$ws.Cells.Item(27, 2).Value2 = 'Monero'
$ws.Cells.Item(27, 3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = '152.35'
$ws.Cells.Item(27, 5).Value2 = '  +1.02%  '
$ws.Cells.Item(28, 2).Value2 = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = '17.89'
$ws.Cells.Item(28, 5).Value2 = '  -2.01%  '
$ws.Cells.Item(29, 2).Value2 = 'LidoDAOToken'
$ws.Cells.Item(29, 3).Value2 = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = '2.031'
$ws.Cells.Item(29, 5).Value2 = '  +1.22%  '
$ws.Cells.Item(30, 2).Value2 = 'BitcoinCash'
$ws.Cells.Item(30, 3).Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = '111.94'
$ws.Cells.Item(30, 5).Value2 = '  -1.45%  '
$ws.Cells.Item(31, 2).Value2 = 'InternetComputer(DFINITY)'
$ws.Cells.Item(31, 3).Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = '4.846'
$ws.Cells.Item(31, 5).Value2 = '  -2.29%  '
$ws.Cells.Item(32, 2).Value2 = 'Stellar'
$ws.Cells.Item(32, 3).Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = '0.08671'
$ws.Cells.Item(32, 5).Value2 = '  -1.95%  '
$ws.Cells.Item(33, 2).Value2 = 'HuobiToken'
$ws.Cells.Item(33, 3).Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = '3.079'
$ws.Cells.Item(33, 5).Value2 = '  -1.66%  '
$ws.Cells.Item(34, 2).Value2 = 'Filecoin'
$ws.Cells.Item(34, 3).Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = '4.443'
$ws.Cells.Item(34, 5).Value2 = '  -0.72%  '
$ws.Cells.Item(35, 2).Value2 = 'ImmutableX'
$ws.Cells.Item(35, 3).Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = '0.7254'
$ws.Cells.Item(35, 5).Value2 = '  -5.63%  '
$ws.Cells.Item(36, 2).Value2 = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = '1.110'
$ws.Cells.Item(36, 5).Value2 = '  -2.25%  '
$ws.Cells.Item(37, 2).Value2 = 'RenderToken'
$ws.Cells.Item(37, 3).Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = '2.655'
$ws.Cells.Item(37, 5).Value2 = '  -0.36%  '
$ws.Cells.Item(38, 2).Value2 = 'Frax'
$ws.Cells.Item(38, 3).Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value2 = '1.006'
$ws.Cells.Item(38, 5).Value2 = '  +0.42%  '
$ws.Cells.Item(39, 2).Value2 = 'TrustWalletToken'
$ws.Cells.Item(39, 3).Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '1.076'
$ws.Cells.Item(39, 5).Value2 = '  -1.61%  '
$ws.Cells.Item(40, 2).Value2 = 'VeChain'
$ws.Cells.Item(40, 3).Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '0.01938'
$ws.Cells.Item(40, 5).Value2 = '  +0.39%  '
$ws.Cells.Item(41, 2).Value2 = 'Hedera'
$ws.Cells.Item(41, 3).Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = '0.05107'
$ws.Cells.Item(41, 5).Value2 = '  -0.63%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = '2.868'
$ws.Cells.Item(42, 5).Value2 = '  -1.99%  '
$ws.Cells.Item(43, 2).Value2 = 'TheSandbox'
$ws.Cells.Item(43, 3).Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = '0.5226'
$ws.Cells.Item(43, 5).Value2 = '  +4.75%  '
$ws.Cells.Item(44, 2).Value2 = 'FraxShare'
$ws.Cells.Item(44, 3).Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = '6.885'
$ws.Cells.Item(44, 5).Value2 = '  -0.98%  '
$ws.Cells.Item(45, 2).Value2 = 'Algorand'
$ws.Cells.Item(45, 3).Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = '0.1522'
$ws.Cells.Item(45, 5).Value2 = '  -5.02%  '
$ws.Cells.Item(46, 2).Value2 = 'Aptos'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '8.022'
$ws.Cells.Item(46, 5).Value2 = '  -4.11%  '
$ws.Cells.Item(47, 2).Value2 = 'Decentraland'
$ws.Cells.Item(47, 3).Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = '0.4702'
$ws.Cells.Item(47, 5).Value2 = '  +0.28%  '
$ws.Cells.Item(48, 2).Value2 = 'PaxDollar'
$ws.Cells.Item(48, 3).Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = '1.008'
$ws.Cells.Item(48, 5).Value2 = '  +0.39%  '
$ws.Cells.Item(49, 2).Value2 = 'EnergySwap'
$ws.Cells.Item(49, 3).Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = '9.893'
$ws.Cells.Item(49, 5).Value2 = '  -3.01%  '
$ws.Cells.Item(50, 2).Value2 = 'Quant'
$ws.Cells.Item(50, 3).Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = '99.91'
$ws.Cells.Item(50, 5).Value2 = '  -2.86%  '
$ws.Cells.Item(51, 2).Value2 = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = '1.584'
$ws.Cells.Item(51, 5).Value2 = '  -2.13%  '
